$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (Hora/hour): every data row (2-51) goes from 21 to 22.
$ws.Range("G2:G51").Value = "'22"

# Column D (Price) and E (Volume(1h)) updates - rows with real numeric data
# (rows 27-37 keep their placeholder "--" / "--%" values and are skipped).
$ws.Range("D2").Value = "'308.41"
$ws.Range("E2").Value = "'-5.85%"
$ws.Range("D3").Value = "'39.94"
$ws.Range("E3").Value = "'-10.27%"
$ws.Range("D4").Value = "'5.059"
$ws.Range("E4").Value = "'-5.92%"
$ws.Range("D5").Value = "'0.07779"
$ws.Range("E5").Value = "'-7.04%"
$ws.Range("D6").Value = "'4.335"
$ws.Range("E6").Value = "'-1.78%"
$ws.Range("D7").Value = "'1.633"
$ws.Range("E7").Value = "'-15.48%"
$ws.Range("D8").Value = "'0.9152"
$ws.Range("E8").Value = "'-5.81%"
$ws.Range("D9").Value = "'0.1005"
$ws.Range("E9").Value = "'-10.62%"
$ws.Range("D10").Value = "'0.1742"
$ws.Range("E10").Value = "'-8.25%"
$ws.Range("D11").Value = "'0.09053"
$ws.Range("E11").Value = "'-6.66%"
$ws.Range("D12").Value = "'0.04456"
$ws.Range("E12").Value = "'-3.29%"
$ws.Range("D13").Value = "'7.039"
$ws.Range("E13").Value = "'-17.41%"
$ws.Range("D14").Value = "'0.1056"
$ws.Range("E14").Value = "'-0.54%"
$ws.Range("D15").Value = "'0.001248"
$ws.Range("E15").Value = "'-4.44%"
$ws.Range("D16").Value = "'0.005764"
$ws.Range("E16").Value = "'-2.11%"
$ws.Range("D17").Value = "'3.353"
$ws.Range("E17").Value = "'-0.42%"
$ws.Range("D18").Value = "'2.562"
$ws.Range("E18").Value = "'1.55%"
$ws.Range("D19").Value = "'0.3369"
$ws.Range("E19").Value = "'0.34%"
$ws.Range("D20").Value = "'0.1388"
$ws.Range("E20").Value = "'0.02%"
$ws.Range("D21").Value = "'0.2659"
$ws.Range("E21").Value = "'3.30%"
$ws.Range("D22").Value = "'0.04150"
$ws.Range("E22").Value = "'-0.08%"
$ws.Range("D23").Value = "'0.001208"
$ws.Range("E23").Value = "'-2.26%"
$ws.Range("D24").Value = "'0.004089"
$ws.Range("E24").Value = "'-7.15%"
$ws.Range("D25").Value = "'0.0001224"
$ws.Range("E25").Value = "'-5.89%"
$ws.Range("D26").Value = "'0.0002989"
$ws.Range("E26").Value = "'0.30%"
$ws.Range("D38").Value = "'0.02411"
$ws.Range("E38").Value = "'-11.19%"
$ws.Range("D39").Value = "'0.05189"
$ws.Range("E39").Value = "'-7.71%"
$ws.Range("D40").Value = "'0.007975"
$ws.Range("E40").Value = "'1.93%"
$ws.Range("D41").Value = "'0.1324"
$ws.Range("E41").Value = "'-6.35%"
$ws.Range("D42").Value = "'0.007106"
$ws.Range("E42").Value = "'-2.89%"
$ws.Range("D43").Value = "'0.002011"
$ws.Range("E43").Value = "'-2.46%"
$ws.Range("D44").Value = "'0.008035"
$ws.Range("E44").Value = "'-7.87%"
$ws.Range("D45").Value = "'0.3340"
$ws.Range("E45").Value = "'-4.82%"
$ws.Range("D46").Value = "'0.00006726"
$ws.Range("E46").Value = "'-2.39%"
$ws.Range("D47").Value = "'0.00000000753"
$ws.Range("E47").Value = "'0.30%"
$ws.Range("D48").Value = "'0.003397"
$ws.Range("E48").Value = "'-2.62%"
$ws.Range("D49").Value = "'0.004115"
$ws.Range("E49").Value = "'16.49%"
$ws.Range("D50").Value = "'0.00002108"
$ws.Range("E50").Value = "'0.30%"
$ws.Range("D51").Value = "'0.0002007"
$ws.Range("E51").Value = "'0.30%"
